$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header cell in A1 (pandas-style "Unnamed: 0" column label)
$ws.Range("A1").Value = "Unnamed: 0"

# Clear the "State" column values (A2:A55). Using a leading apostrophe keeps
# the cells typed as text (matching the original inline-string cell type)
# instead of turning them into empty/blank numeric cells, then the style is
# reset to drop the quote-prefix formatting that the apostrophe trick adds.
$rng = $ws.Range("A2:A55")
$rng.Value = "'"
$rng.Style = "Normal"
